# Update Name of Algo
# Apply updated numeric values (columns A and D) produced by the new
# algorithm run, on Sheet1 of the active workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D4").Value = -8.016
$ws.Range("A9").Value = -21.847
$ws.Range("D9").Value = -8.123000000000001
$ws.Range("D11").Value = -7.323
$ws.Range("A13").Value = -22.08
$ws.Range("A16").Value = -21.905
$ws.Range("D16").Value = -8.493
$ws.Range("A18").Value = -22.091
$ws.Range("A20").Value = -20.654
$ws.Range("D23").Value = -8.348000000000001
$ws.Range("D24").Value = -7.170999999999999
$ws.Range("A26").Value = -21.303
$ws.Range("D26").Value = -7.543000000000001
$ws.Range("A27").Value = -21.874
$ws.Range("A29").Value = -21.393
$ws.Range("D34").Value = -7.636
$ws.Range("A35").Value = -20.089
$ws.Range("D35").Value = -8.254000000000001
$ws.Range("A36").Value = -20.391
$ws.Range("D44").Value = -7.88
$ws.Range("A45").Value = -21.576
$ws.Range("D48").Value = -7.986
$ws.Range("D49").Value = -8.35
$ws.Range("D52").Value = -7.958
$ws.Range("A55").Value = -22.003
$ws.Range("A57").Value = -22.232
$ws.Range("D66").Value = -7.486
$ws.Range("D67").Value = -7.656000000000001
$ws.Range("A69").Value = -21.52
$ws.Range("D73").Value = -8.193999999999999
$ws.Range("A76").Value = -20.074
$ws.Range("A78").Value = -19.911
$ws.Range("D78").Value = -8.379999999999999
$ws.Range("D80").Value = -8.044
$ws.Range("A82").Value = -21.901
$ws.Range("A83").Value = -21.891
$ws.Range("D91").Value = -7.376
$ws.Range("A93").Value = -21.461
$ws.Range("A97").Value = -21.349
$ws.Range("D97").Value = -7.543000000000001
$ws.Range("D99").Value = -8.253
$ws.Range("D104").Value = -7.790999999999999
